# This script applies numeric corrections to the "想去人数" (F column) and
# "最低票价" (G column) figures across the four worksheets of the workbook,
# matching the upstream data refresh captured in the commit
# "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 3413
$ws.Range("F4").Value = 578
$ws.Range("G4").Value = 70
$ws.Range("F5").Value = 835
$ws.Range("G5").Value = 70
$ws.Range("F6").Value = 316
$ws.Range("G6").Value = 60
$ws.Range("F7").Value = 273
$ws.Range("F8").Value = 74
$ws.Range("F9").Value = 161
$ws.Range("F10").Value = 631
$ws.Range("F12").Value = 423
$ws.Range("F13").Value = 66
$ws.Range("F14").Value = 491
$ws.Range("F15").Value = 326
$ws.Range("F19").Value = 183

# --- Sheet 2: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 6
$ws.Range("F16").Value = 25

# --- Sheet 3: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 6204
$ws.Range("F5").Value = 1780
$ws.Range("F6").Value = 114

# --- Sheet 4: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 6204
$ws.Range("F5").Value = 1780
$ws.Range("F6").Value = 3413
$ws.Range("F7").Value = 114
$ws.Range("F9").Value = 578
$ws.Range("G9").Value = 70
$ws.Range("F10").Value = 835
$ws.Range("G10").Value = 70
$ws.Range("F11").Value = 316
$ws.Range("G11").Value = 60
$ws.Range("F12").Value = 273
$ws.Range("F14").Value = 74
$ws.Range("F17").Value = 161
$ws.Range("F19").Value = 6
$ws.Range("F20").Value = 631
$ws.Range("F24").Value = 423
$ws.Range("F26").Value = 66
$ws.Range("F27").Value = 491
$ws.Range("F29").Value = 326
$ws.Range("F38").Value = 25
$ws.Range("F40").Value = 183
